$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.97750833333333
$ws.Range("H2").Value = 65.932525
$ws.Range("I2").Value = 0.5427578249542736
$ws.Range("J2").Value = 0.5427578249542736
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.854571666666667
$ws.Range("N2").Value = 5.563715
$ws.Range("O2").Value = 0.01651371646154392
$ws.Range("P2").Value = 0.01651371646154392
$ws.Range("Q2").Value = 40.75886425893056
$ws.Range("R2").Value = 366.829778330375
$ws.Range("S2").Value = 0.00896294882857916
$ws.Range("T2").Value = 0.008962948828579162

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.97750833333333
$ws.Range("H3").Value = 65.932525
$ws.Range("I3").Value = 0.5427578249542736
$ws.Range("J3").Value = 0.5427578249542736
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 92.91372433333333
$ws.Range("N3").Value = 278.741173
$ws.Range("O3").Value = 0.8273343794712995
$ws.Range("P3").Value = 0.8273343794712996
$ws.Range("Q3").Value = 2042.01215081687
$ws.Range("R3").Value = 18378.10935735183
$ws.Range("S3").Value = 0.4490422083117361
$ws.Range("T3").Value = 0.4490422083117362

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.97750833333333
$ws.Range("H4").Value = 65.932525
$ws.Range("I4").Value = 0.5427578249542736
$ws.Range("J4").Value = 0.5427578249542736
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.162136
$ws.Range("N4").Value = 0.4864080000000001
$ws.Range("O4").Value = 0.001443712303133186
$ws.Range("P4").Value = 0.001443712303133187
$ws.Range("Q4").Value = 3.563345291133334
$ws.Range("R4").Value = 32.07010762020001
$ws.Range("S4").Value = 0.0007835861495082931
$ws.Range("T4").Value = 0.0007835861495082932

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.97750833333333
$ws.Range("H5").Value = 65.932525
$ws.Range("I5").Value = 0.5427578249542736
$ws.Range("J5").Value = 0.5427578249542736
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.37449166666667
$ws.Range("N5").Value = 52.123475
$ws.Range("O5").Value = 0.1547081917640233
$ws.Range("P5").Value = 0.1547081917640233
$ws.Range("Q5").Value = 381.8480353915973
$ws.Range("R5").Value = 3436.632318524375
$ws.Range("S5").Value = 0.08396908166444995
$ws.Range("T5").Value = 0.08396908166444995

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 12.688376
$ws.Range("H6").Value = 38.065128
$ws.Range("I6").Value = 0.3133528721960219
$ws.Range("J6").Value = 0.3133528721960219
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.854571666666667
$ws.Range("N6").Value = 5.563715
$ws.Range("O6").Value = 0.01651371646154392
$ws.Range("P6").Value = 0.01651371646154392
$ws.Range("Q6").Value = 23.53150262561334
$ws.Range("R6").Value = 211.78352363052
$ws.Range("S6").Value = 0.005174620483855515
$ws.Range("T6").Value = 0.005174620483855516

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 12.688376
$ws.Range("H7").Value = 38.065128
$ws.Range("I7").Value = 0.3133528721960219
$ws.Range("J7").Value = 0.3133528721960219
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 92.91372433333333
$ws.Range("N7").Value = 278.741173
$ws.Range("O7").Value = 0.8273343794712995
$ws.Range("P7").Value = 0.8273343794712996
$ws.Range("Q7").Value = 1178.924269901683
$ws.Range("R7").Value = 10610.31842911514
$ws.Range("S7").Value = 0.2592476040738452
$ws.Range("T7").Value = 0.2592476040738452

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.688376
$ws.Range("H8").Value = 38.065128
$ws.Range("I8").Value = 0.3133528721960219
$ws.Range("J8").Value = 0.3133528721960219
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.162136
$ws.Range("N8").Value = 0.4864080000000001
$ws.Range("O8").Value = 0.001443712303133186
$ws.Range("P8").Value = 0.001443712303133187
$ws.Range("Q8").Value = 2.057242531136001
$ws.Range("R8").Value = 18.515182780224
$ws.Range("S8").Value = 0.0004523913968115178
$ws.Range("T8").Value = 0.0004523913968115178

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.688376
$ws.Range("H9").Value = 38.065128
$ws.Range("I9").Value = 0.3133528721960219
$ws.Range("J9").Value = 0.3133528721960219
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 17.37449166666667
$ws.Range("N9").Value = 52.123475
$ws.Range("O9").Value = 0.1547081917640233
$ws.Range("P9").Value = 0.1547081917640233
$ws.Range("Q9").Value = 220.4540830755333
$ws.Range("R9").Value = 1984.0867476798
$ws.Range("S9").Value = 0.04847825624150965
$ws.Range("T9").Value = 0.04847825624150965

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.008189666666667
$ws.Range("H10").Value = 6.024569
$ws.Range("I10").Value = 0.04959436889042158
$ws.Range("J10").Value = 0.04959436889042158
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.854571666666667
$ws.Range("N10").Value = 5.563715
$ws.Range("O10").Value = 0.01651371646154392
$ws.Range("P10").Value = 0.01651371646154392
$ws.Range("Q10").Value = 3.724331657092778
$ws.Range("R10").Value = 33.518984913835
$ws.Range("S10").Value = 0.0008189873459456364
$ws.Range("T10").Value = 0.0008189873459456366

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.008189666666667
$ws.Range("H11").Value = 6.024569
$ws.Range("I11").Value = 0.04959436889042158
$ws.Range("J11").Value = 0.04959436889042158
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 92.91372433333333
$ws.Range("N11").Value = 278.741173
$ws.Range("O11").Value = 0.8273343794712995
$ws.Range("P11").Value = 0.8273343794712996
$ws.Range("Q11").Value = 186.5883810977152
$ws.Range("R11").Value = 1679.295429879437
$ws.Range("S11").Value = 0.04103112641122766
$ws.Range("T11").Value = 0.04103112641122766

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.008189666666667
$ws.Range("H12").Value = 6.024569
$ws.Range("I12").Value = 0.04959436889042158
$ws.Range("J12").Value = 0.04959436889042158
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.162136
$ws.Range("N12").Value = 0.4864080000000001
$ws.Range("O12").Value = 0.001443712303133186
$ws.Range("P12").Value = 0.001443712303133187
$ws.Range("Q12").Value = 0.3255998397946667
$ws.Range("R12").Value = 2.930398558152
$ws.Range("S12").Value = [double]"7.160000053322739E-05"
$ws.Range("T12").Value = [double]"7.160000053322739E-05"

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.008189666666667
$ws.Range("H13").Value = 6.024569
$ws.Range("I13").Value = 0.04959436889042158
$ws.Range("J13").Value = 0.04959436889042158
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 17.37449166666667
$ws.Range("N13").Value = 52.123475
$ws.Range("O13").Value = 0.1547081917640233
$ws.Range("P13").Value = 0.1547081917640233
$ws.Range("Q13").Value = 34.89127462858611
$ws.Range("R13").Value = 314.021471657275
$ws.Range("S13").Value = 0.007672655132715053
$ws.Range("T13").Value = 0.007672655132715053

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.818218
$ws.Range("H14").Value = 11.454654
$ws.Range("I14").Value = 0.09429493395928291
$ws.Range("J14").Value = 0.09429493395928291
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.854571666666667
$ws.Range("N14").Value = 5.563715
$ws.Range("O14").Value = 0.01651371646154392
$ws.Range("P14").Value = 0.01651371646154392
$ws.Range("Q14").Value = 7.081158919956667
$ws.Range("R14").Value = 63.73043027961
$ws.Range("S14").Value = 0.001557159803163607
$ws.Range("T14").Value = 0.001557159803163607

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.818218
$ws.Range("H15").Value = 11.454654
$ws.Range("I15").Value = 0.09429493395928291
$ws.Range("J15").Value = 0.09429493395928291
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 92.91372433333333
$ws.Range("N15").Value = 278.741173
$ws.Range("O15").Value = 0.8273343794712995
$ws.Range("P15").Value = 0.8273343794712996
$ws.Range("Q15").Value = 354.7648546965713
$ws.Range("R15").Value = 3192.883692269142
$ws.Range("S15").Value = 0.07801344067449049
$ws.Range("T15").Value = 0.07801344067449051

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.818218
$ws.Range("H16").Value = 11.454654
$ws.Range("I16").Value = 0.09429493395928291
$ws.Range("J16").Value = 0.09429493395928291
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.162136
$ws.Range("N16").Value = 0.4864080000000001
$ws.Range("O16").Value = 0.001443712303133186
$ws.Range("P16").Value = 0.001443712303133187
$ws.Range("Q16").Value = 0.6190705936480001
$ws.Range("R16").Value = 5.571635342832001
$ws.Range("S16").Value = 0.000136134756280148
$ws.Range("T16").Value = 0.0001361347562801481

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.818218
$ws.Range("H17").Value = 11.454654
$ws.Range("I17").Value = 0.09429493395928291
$ws.Range("J17").Value = 0.09429493395928291
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.37449166666667
$ws.Range("N17").Value = 52.123475
$ws.Range("O17").Value = 0.1547081917640233
$ws.Range("P17").Value = 0.1547081917640233
$ws.Range("Q17").Value = 66.33959682251667
$ws.Range("R17").Value = 597.05637140265
$ws.Range("S17").Value = 0.01458819872534865
$ws.Range("T17").Value = 0.01458819872534865

Write-Host "Applied all changes"